# Apply odds updates for rows 3 and 4 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value  = 1.25
$ws.Range("H3").Value  = 5.75
$ws.Range("I3").Value  = 11
$ws.Range("J3").Value  = 1.67
$ws.Range("L3").Value  = 8
$ws.Range("Y3").Value  = 9.5
$ws.Range("AD3").Value = 11
$ws.Range("AG3").Value = 29
$ws.Range("AI3").Value = 29
$ws.Range("AJ3").Value = 126
$ws.Range("AU3").Value = 9
$ws.Range("AW3").Value = 10
$ws.Range("AZ3").Value = 151

# Row 4 updates
$ws.Range("I4").Value  = 5.25
$ws.Range("L4").Value  = 5.5
$ws.Range("M4").Value  = 1.1
$ws.Range("N4").Value  = 7
$ws.Range("AI4").Value = 17
$ws.Range("AQ4").Value = 34
$ws.Range("BB4").Value = 451
